# Code_Structure.xlsx update:
# Reviewed two scripts (Inputs.py and config.json) and filled in the
# "basic purpose", "non-obvious details" and "suggestions" columns for
# those rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$D2 = "To read in parameters, primarily from config.json (but default parameters can be read from config_Edsizes.json and config_Edtypes.json if not supplied)"

$E2 = @"
Functions are as follows:
parse_config_acuity - Checks for errors on acuity percentages in config file
parse_config_simtype - Checks if sim_parameters (daily number of patients, numbers of clinicans, maximum number of seats) has errors in config file
parse_config_Sim_Parameters - Checks if several other sim_parameters have errors in config file
get_EDtype_parameters - Return acuity percentages and various other sim parameters. 
get_EDsize_parameters - Return ED size parameters
get_target_time - Return target time, which affects time for decision to admit
parse_bed_occupancy - Get bed time, which is time after decision to admit to actually get a bed.
get_doctor_times - Get doctor times 
get_diagnostic_percentages - Get diagnostic percentages
-read_parameters - Read in input data using above functions.
"@

$D3 = "To provide model parameters"

$E3 = @"
Max_Waiting_Times - Max waiting time of patient before leaving  - i.e. their patience
Prob_Admission - Is probability is admitted or dicharged
"@

$F3 = @"
Do we want to have patient patience?
Opening hours on rooms (SDEC has these).
Numbers of doctors available doesn't depend on time of day
Might want surge capacity beds because strict capacity constraints on certain rooms is likely unrealistic
Many parameters are largely static (do not change depending on time of day) - need to think what we would want to model
"@

$ws.Range("D2").Value = $D2
$ws.Range("D3").Value = $D3
$ws.Range("E2").Value = $E2
$ws.Range("E3").Value = $E3
$ws.Range("F3").Value = $F3

# Wrap text for the newly-filled long comment cells (D3 keeps the default style).
$ws.Range("D2").WrapText = $true
$ws.Range("E2").WrapText = $true
$ws.Range("E3").WrapText = $true
$ws.Range("F3").WrapText = $true

# Row heights grow to accommodate the wrapped multi-line text.
$ws.Rows.Item(2).RowHeight = 302.4
$ws.Rows.Item(3).RowHeight = 100.8

# Column E widens to fit the new comments (no longer auto "best fit").
$ws.Columns.Item(5).ColumnWidth = 41.83

# Selection ends on F4, matching where the reviewer finished editing.
[void]$ws.Range("F4").Select()
